$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.617.52"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "2.034.54"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'234.14"
$ws.Range("E5").Value = "  -9.28%  "
$ws.Range("D6").Value = "'0.601"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'55.29"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").Value = "'57.74"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("D11").Value = "'0.0752"
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "2.333.79"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("E15").Value = "  -4.76%  "
$ws.Range("D16").Value = "'0.764"
$ws.Range("E16").Value = "  -3.67%  "
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "2.034.86"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("D19").Value = "36.760.56"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "'67.83"
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").Value = "'5.49"
$ws.Range("E21").Value = "  +8.19%  "
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "'221.34"
$ws.Range("E23").Value = "  -5.46%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  -5.80%  "
$ws.Range("D27").Value = "'163.22"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("E28").Value = "  +8.01%  "
$ws.Range("D29").Value = "'8.66"
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("D30").Value = "'19.03"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("E34").Value = "  -5.28%  "
$ws.Range("D35").Value = "'2.46"
$ws.Range("E35").Value = "  +5.11%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'3.31"
$ws.Range("E39").Value = "  -4.88%  "
$ws.Range("D40").Value = "'5.81"
$ws.Range("E40").Value = "  +5.06%  "
$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").Value = "'0.0945"
$ws.Range("E42").Value = "  +3.03%  "
$ws.Range("D43").Value = "1.459.27"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0204"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.14"
$ws.Range("E45").Value = "  +39.40%  "
$ws.Range("D46").Value = "'91.04"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").Value = "'1.12"
$ws.Range("E47").Value = "  -4.82%  "
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "'2.90"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").Value = "'6.88"
$ws.Range("E51").Value = "  -0.35%  "
